# Applies the "new changes for poc" edit to bluenile.xlsx:
#   1. TC001 rows (2-7): "Enabled" column (B) flips from "no" to "yes"
#   2. TC002 row 8: the goto target in column D changes from the generic
#      Amazon homepage to a specific product page URL
#   3. Leaves the active selection on D8 (matches the saved UI state)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Flip "Enabled" to "yes" for the TC001 block (rows 2-7, column B)
$ws.Range("B2:B7").Value = "yes"

# 2) Update the TC002 goto URL (row 8, column D) to the new product page
$ws.Range("D8").Value = "https://www.amazon.com/SIMPLIHOME-Cocktail-Footrest-Upholstered-Contemporary/dp/B06WV91PNR?th=1"

# 3) Match the saved cursor/selection position from the edit
$ws.Range("D8").Select()
